# Update forecast dates/values in "Forecast Comparison" (shift week-start dates
# forward by one week) and refresh dependent figures in "Summary".

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to stay a text value (avoid Excel auto-parsing things
    # like "2025-02-02" into a date serial number), then restore the
    # original "Normal" style so no stray number format sticks around.
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

# ---- "Forecast Comparison" sheet: new Week_Start_Date values (row -> date) ----
$newDates = @{
    2  = "2025-02-02"
    3  = "2025-02-09"
    4  = "2025-02-16"
    5  = "2025-02-23"
    6  = "2025-03-02"
    7  = "2025-03-09"
    8  = "2025-03-16"
    9  = "2025-03-23"
    10 = "2025-03-30"
    11 = "2025-04-06"
    12 = "2025-04-13"
    13 = "2025-04-20"
    14 = "2025-04-27"
    15 = "2025-05-04"
    16 = "2025-05-11"
    17 = "2025-05-18"
}

foreach ($row in 2..17) {
    Set-TextValue ($wsForecast.Cells.Item($row, 2)) $newDates[$row]
}

# ---- "Forecast Comparison" sheet: numeric forecast value corrections ----
$wsForecast.Range("E9").Value2  = 2
$wsForecast.Range("H9").Value2  = 8

$wsForecast.Range("F11").Value2 = 2

$wsForecast.Range("D12").Value2 = 0
$wsForecast.Range("E12").Value2 = 3
$wsForecast.Range("F12").Value2 = 2
$wsForecast.Range("H12").Value2 = 11

$wsForecast.Range("H13").Value2 = 12

$wsForecast.Range("H14").Value2 = 12

$wsForecast.Range("G15").Value2 = 5
$wsForecast.Range("H15").Value2 = 11

$wsForecast.Range("G16").Value2 = 6
$wsForecast.Range("H16").Value2 = 13

$wsForecast.Range("H17").Value2 = 12

# ---- "Summary" sheet updates ----
Set-TextValue ($wsSummary.Range("B2"))  "2022-12-25 to 2025-01-26"
Set-TextValue ($wsSummary.Range("B13")) "2025-03-16"
Set-TextValue ($wsSummary.Range("B15")) "2025-02-02"
